# Actualización automática 2025-11-17 16:30:09
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("M12").Value = 1466.07
$wsVentasGrupo.Range("I31").Value = 221
$wsVentasGrupo.Range("M49").Value = 1588.43
$wsVentasGrupo.Range("I60").Value = "3 de 58"
$wsVentasGrupo.Range("M60").Value = "8 de 58"

# --- Sheet "VENTA MENSUAL" ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F12").Value = 1466.07
$wsVentaMensual.Range("F31").Value = 1621.66
$wsVentaMensual.Range("F49").Value = 1588.43
$wsVentaMensual.Range("F60").Value = 15772.57

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D7").Value = 456.2
$wsCumplimiento.Range("E7").Value = 430.511016287574
$wsCumplimiento.Range("F7").Value = 0.5144855444674518

$wsCumplimiento.Range("D12").Value = 7161.18
$wsCumplimiento.Range("E12").Value = 43145.82
$wsCumplimiento.Range("F12").Value = 0.1423495736179856

$wsCumplimiento.Range("D14").Value = 17344.7
$wsCumplimiento.Range("E14").Value = 80517.18766749099
$wsCumplimiento.Range("F14").Value = 0.1772365158020734
